# Auto-generated PowerShell COM-interop script
# Applies the cached-value updates described by the commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Table_*" sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 6520.3706
$ws.Range("I15").Value = 6520.3706
$ws.Range("K15").Value = 19561.1118
$ws.Range("M15").Value = -19392.1118

$ws.Range("H57").Value = 50124
$ws.Range("J57").Value = 50124
$ws.Range("L57").Value = 150372
$ws.Range("N57").Value = -151370

$ws.Range("H58").Value = 6495669
$ws.Range("I58").Value = 35714492
$ws.Range("J58").Value = 2597.111
$ws.Range("K58").Value = 107143476
$ws.Range("L58").Value = 7791.333
$ws.Range("M58").Value = -107143326
$ws.Range("N58").Value = -8091.333

$ws.Range("H96").Value = 6592.1377
$ws.Range("I96").Value = 1892.5652
$ws.Range("J96").Value = 24607.166
$ws.Range("K96").Value = 5677.6956
$ws.Range("L96").Value = 73821.49800000001
$ws.Range("M96").Value = -4304.6956
$ws.Range("N96").Value = -76567.49800000001

$ws.Range("H107").Value = 6759.5264
$ws.Range("I107").Value = 7001.722
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 7001.722
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -5081.722
$ws.Range("N107").Value = -6240

$ws.Range("H137").Value = 2140742.5
$ws.Range("I137").Value = 4049490.8
$ws.Range("J137").Value = 7435.4116
$ws.Range("K137").Value = 12148472.4
$ws.Range("L137").Value = 22306.2348
$ws.Range("M137").Value = -12145922.4
$ws.Range("N137").Value = -27406.2348

$ws.Range("H141").Value = 4656
$ws.Range("I141").Value = 2963
$ws.Range("J141").Value = 8042
$ws.Range("K141").Value = 8889
$ws.Range("L141").Value = 24126
$ws.Range("M141").Value = -3709
$ws.Range("N141").Value = -34486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1871.862
$ws.Range("I2").Value = 2012.3182
$ws.Range("J2").Value = 1430.4286
$ws.Range("K2").Value = 2012.3182
$ws.Range("L2").Value = 1430.4286
$ws.Range("M2").Value = -1899.3182
$ws.Range("N2").Value = -1656.4286

$ws.Range("H31").Value = 10428.363
$ws.Range("I31").Value = 2452
$ws.Range("K31").Value = 2452
$ws.Range("M31").Value = -2158

$ws.Range("H32").Value = 11370.777
$ws.Range("I32").Value = 9745.254999999999
$ws.Range("J32").Value = 39004.668
$ws.Range("K32").Value = 9745.254999999999
$ws.Range("L32").Value = 39004.668
$ws.Range("M32").Value = -9458.254999999999
$ws.Range("N32").Value = -39578.668

$ws.Range("H61").Value = 2685.7715
$ws.Range("I61").Value = 1942.5294
$ws.Range("J61").Value = 3387.7222
$ws.Range("K61").Value = 1942.5294
$ws.Range("L61").Value = 3387.7222
$ws.Range("M61").Value = -1730.5294
$ws.Range("N61").Value = -3811.7222

$ws.Range("H116").Value = 1871.862
$ws.Range("I116").Value = 2012.3182
$ws.Range("J116").Value = 1430.4286
$ws.Range("K116").Value = 2012.3182
$ws.Range("L116").Value = 1430.4286
$ws.Range("M116").Value = 281.6818000000001
$ws.Range("N116").Value = -6018.4286

$ws.Range("H132").Value = 9617541
$ws.Range("I132").Value = 13514902
$ws.Range("J132").Value = 4050.6667
$ws.Range("K132").Value = 40544706
$ws.Range("L132").Value = 12152.0001
$ws.Range("M132").Value = -40542176
$ws.Range("N132").Value = -17212.0001

$ws.Range("H136").Value = 2685.7715
$ws.Range("I136").Value = 1942.5294
$ws.Range("J136").Value = 3387.7222
$ws.Range("K136").Value = 5827.5882
$ws.Range("L136").Value = 10163.1666
$ws.Range("M136").Value = -3277.5882
$ws.Range("N136").Value = -15263.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1871.862
$ws.Range("I3").Value = 2012.3182
$ws.Range("J3").Value = 1430.4286
$ws.Range("K3").Value = 2012.3182
$ws.Range("L3").Value = 1430.4286
$ws.Range("M3").Value = -1898.3182
$ws.Range("N3").Value = -1658.4286

$ws.Range("H82").Value = 20000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 20000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -20766

$ws.Range("H85").Value = 20000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 20000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -22652

$ws.Range("H134").Value = 3246.6
$ws.Range("I134").Value = 2441.2856
$ws.Range("J134").Value = 3680.2307
$ws.Range("K134").Value = 7323.8568
$ws.Range("L134").Value = 11040.6921
$ws.Range("M134").Value = -4788.8568
$ws.Range("N134").Value = -16110.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9265813
$ws.Range("I31").Value = 2605.5789
$ws.Range("J31").Value = 19618810
$ws.Range("K31").Value = 2605.5789
$ws.Range("L31").Value = 19618810
$ws.Range("M31").Value = -2310.5789
$ws.Range("N31").Value = -19619400

$ws.Range("H34").Value = 9265813
$ws.Range("I34").Value = 2605.5789
$ws.Range("J34").Value = 19618810
$ws.Range("K34").Value = 2605.5789
$ws.Range("L34").Value = 19618810
$ws.Range("M34").Value = -2403.5789
$ws.Range("N34").Value = -19619214

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H105").Value = 3685.7144
$ws.Range("I105").Value = 3800
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 3800
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -2053
$ws.Range("N105").Value = -6494

$ws.Range("H132").Value = 96310.734
$ws.Range("I132").Value = 2136
$ws.Range("J132").Value = 178713.62
$ws.Range("K132").Value = 6408
$ws.Range("L132").Value = 536140.86
$ws.Range("M132").Value = -3878
$ws.Range("N132").Value = -541200.86

$ws.Range("H139").Value = 27778
$ws.Range("I139").Value = 19900
$ws.Range("J139").Value = 29747.5
$ws.Range("K139").Value = 19900
$ws.Range("L139").Value = 29747.5
$ws.Range("M139").Value = -14760
$ws.Range("N139").Value = -40027.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5378.7144

$ws.Range("H131").Value = 2397.305
$ws.Range("I131").Value = 14718.429
$ws.Range("J131").Value = 1247.3334
$ws.Range("K131").Value = 44155.287
$ws.Range("L131").Value = 3742.0002
$ws.Range("M131").Value = -39115.287
$ws.Range("N131").Value = -13822.0002

$ws.Range("H135").Value = 5378.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5274.9644
$ws.Range("I70").Value = 5380.905
$ws.Range("J70").Value = 4957.143
$ws.Range("K70").Value = 5380.905
$ws.Range("L70").Value = 4957.143
$ws.Range("M70").Value = -5110.905
$ws.Range("N70").Value = -5497.143

$ws.Range("H73").Value = 5274.9644
$ws.Range("I73").Value = 5380.905
$ws.Range("J73").Value = 4957.143
$ws.Range("K73").Value = 5380.905
$ws.Range("L73").Value = 4957.143
$ws.Range("M73").Value = -4444.905
$ws.Range("N73").Value = -6829.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2859.4614
$ws.Range("I7").Value = 1896
$ws.Range("J7").Value = 4401
$ws.Range("K7").Value = 1896
$ws.Range("L7").Value = 4401
$ws.Range("M7").Value = -1784
$ws.Range("N7").Value = -4625

$ws.Range("H126").Value = 2859.4614
$ws.Range("I126").Value = 1896
$ws.Range("J126").Value = 4401
$ws.Range("K126").Value = 5688
$ws.Range("L126").Value = 13203
$ws.Range("M126").Value = -3218
$ws.Range("N126").Value = -18143

$ws.Range("H138").Value = 58000
$ws.Range("J138").Value = 58000
$ws.Range("L138").Value = 58000
$ws.Range("N138").Value = -68280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5001253
$ws.Range("I107").Value = 1043.3846
$ws.Range("J107").Value = 14287357
$ws.Range("K107").Value = 3130.1538
$ws.Range("L107").Value = 42862071
$ws.Range("M107").Value = -1210.1538
$ws.Range("N107").Value = -42865911

$ws.Range("H132").Value = 1893069.9
$ws.Range("I132").Value = 3346555.2
$ws.Range("J132").Value = 3539
$ws.Range("K132").Value = 10039665.6
$ws.Range("L132").Value = 10617
$ws.Range("M132").Value = -10037135.6
$ws.Range("N132").Value = -15677

$ws.Range("H136").Value = 494879.2
$ws.Range("I136").Value = 667467.5600000001
$ws.Range("J136").Value = 1769.5714
$ws.Range("K136").Value = 2002402.68
$ws.Range("L136").Value = 5308.7142
$ws.Range("M136").Value = -1999852.68
$ws.Range("N136").Value = -10408.7142

